$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Y2").Value = -3.992934669072811
$ws.Range("Z2").Value = 2.042399836814456
$ws.Range("F3").Value = -0.7818945361920984
$ws.Range("G3").Value = -0.3631425419165785
$ws.Range("H3").Value = -0.4445839688494094
$ws.Range("I3").Value = -0.1248950033487057
$ws.Range("M3").Value = -0.5545191045891807
$ws.Range("N3").Value = -0.2025438469602849
$ws.Range("Q3").Value = -0.4445839688494094
$ws.Range("R3").Value = -0.1248950033487057
$ws.Range("S3").Value = -0.7818945361920984
$ws.Range("T3").Value = -0.3631425419165785
$ws.Range("U3").Value = 0.00947832987223423
$ws.Range("V3").Value = 0.01601584133781592
$ws.Range("Y3").Value = -1.026492312377353
$ws.Range("Z3").Value = 0.557779929346007
$ws.Range("AA3").Value = 0.009200832113358808
$ws.Range("AB3").Value = 0.01516755802386971
$ws.Range("AC3").Value = 0.01517041424497499
$ws.Range("AD3").Value = 0.009201766465057748
$ws.Range("F4").Value = -0.6353125423309869
$ws.Range("G4").Value = -0.2312119145187517
$ws.Range("H4").Value = -0.3542204028787461
$ws.Range("I4").Value = -0.0326722990455246
$ws.Range("M4").Value = -0.3899191045891807
$ws.Range("N4").Value = -0.05788683021062446
$ws.Range("Q4").Value = -0.3542204028787461
$ws.Range("R4").Value = -0.0326722990455246
$ws.Range("S4").Value = -0.6353125423309869
$ws.Range("T4").Value = -0.2312119145187517
$ws.Range("U4").Value = 0.01195545863112958
$ws.Range("V4").Value = 0.0201929877776081
$ws.Range("AA4").Value = 0.01147413427983304
$ws.Range("AB4").Value = 0.019314879002934
$ws.Range("AC4").Value = 0.01932024166328655
$ws.Range("AD4").Value = 0.01147707375508925
$ws.Range("AA5").Value = 0.009589471774614205
$ws.Range("AB5").Value = 0.01892729117407374
$ws.Range("AC5").Value = 0.01891359958648576
$ws.Range("AD5").Value = 0.009522895840483667
$ws.Range("AA6").Value = 0.008960160494083536
$ws.Range("AB6").Value = 0.0189033048337237
$ws.Range("AC6").Value = 0.0188987815260182
$ws.Range("AD6").Value = 0.008841485926690635
$ws.Range("AA7").Value = 0.008047744974150978
$ws.Range("AB7").Value = 0.01921647436189348
$ws.Range("AC7").Value = 0.01933476268794958
$ws.Range("AD7").Value = 0.007703113409416113
$ws.Range("Y8").Value = -0.01864904827323405
$ws.Range("Z8").Value = 0.04295436866737521
$ws.Range("AA8").Value = 0.007866244696811138
$ws.Range("AB8").Value = 0.01961008072833311
$ws.Range("AC8").Value = 0.01990403200445929
$ws.Range("AD8").Value = 0.007730697503751485
$ws.Range("F9").Value = -0.3613
$ws.Range("H9").Value = -0.3613
$ws.Range("M9").Value = -0.3613
$ws.Range("Q9").Value = -0.3613
$ws.Range("S9").Value = -0.3613
$ws.Range("U9").Value = 0.002448386286034203
$ws.Range("V9").Value = 0.03240565138774108
$ws.Range("AA9").Value = 0.007697443815743082
$ws.Range("AB9").Value = 0.02000874445065931
$ws.Range("AC9").Value = 0.0204723371638521
$ws.Range("AD9").Value = 0.007501150915082201
